$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entregables")

# Mark the "Agregar", "Eliminar", "Consultar" rows as completed in column C
$ws.Range("C5").Value = "x"
$ws.Range("C6").Value = "x"
$ws.Range("C7").Value = "x"

# Update the active cell selection to reflect where the user left off
$ws.Activate()
$ws.Range("C9").Select()
